{"js": "// 1. Remove the \"Meta description: ...\" paragraph that follows the title (H1).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0 is the H1 title; paragraph 1 is the \"Meta description\" paragraph.\nparagraphs.items[1].delete();\nawait context.sync();\n\n// 2. Insert a new bold paragraph with the title text right before the final\n//    (Prompt / image-description) paragraph.\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs2.items[paragraphs2.items.length - 1];\nlastParagraph.insertHtml(\n  \"<p><b>Play Fruits in Flames Slot Game for Free - Review</b></p>\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// 3. Replace the text of the final paragraph (still italic) with the new\n//    meta-description copy.\nconst paragraphs3 = context.document.body.paragraphs;\nparagraphs3.load(\"items\");\nawait context.sync();\n\nconst finalParagraph = paragraphs3.items[paragraphs3.items.length - 1];\nconst finalRange = finalParagraph.getRange();\nfinalRange.insertText(\n  \"Read our review of Fruits in Flames, a classic slot game with traditional gameplay, high winning potential, and a gamble feature. Play for free.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n\n# 1. Remove the \"Meta description: ...\" paragraph that follows the H1 title.\n$doc.Paragraphs.Item(2).Range.Delete()\n\n# 2. Insert a new bold paragraph with the title text right before the final\n#    (Prompt / image-description) paragraph. Anchor the insertion at the end\n#    of the paragraph just before it, so no paragraph style / run formatting\n#    is inherited from the italic \"Prompt\" paragraph.\n$n = $doc.Paragraphs.Count\n$anchorPara = $doc.Paragraphs.Item($n - 1)\n$insertRange = $anchorPara.Range.Duplicate\n$insertRange.Collapse(0)  # wdCollapseEnd\n$insertRange.InsertAfter(\"Play Fruits in Flames Slot Game for Free - Review`r\")\n$insertRange.Font.Bold = 1\n\n# 3. Replace the text of the final paragraph (still italic) with the new\n#    meta-description copy. Shrink the range by one character first so the\n#    trailing paragraph mark is excluded from the replace, which leaves the\n#    paragraph's leading empty run and italic run formatting untouched.\n$n = $doc.Paragraphs.Count\n$lastPara = $doc.Paragraphs.Item($n)\n$textRange = $lastPara.Range.Duplicate\n$textRange.End = $textRange.End - 1\n$textRange.Text = \"Read our review of Fruits in Flames, a classic slot game with traditional gameplay, high winning potential, and a gamble feature. Play for free.\"\n"}
